$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Both sheets store lat/long coordinates. A bunch of the latitude values
# (column B) lost their decimal point at some point (e.g. 71539 instead of
# 7.1539) -- restore the missing "." after the leading digit. One longitude
# cell (CAI!C4) has the same problem (-73133 instead of -73.133).
# While we're at it, give every column-B cell a plain "General" number
# format so the fixed decimals display correctly instead of as thousands.
# ---------------------------------------------------------------------------

$wsCAI = $wb.Worksheets.Item("CAI")

$caiLat = @{
    2  = 7.1539
    3  = 7.1229
    4  = 7.1237
    5  = 7.1415
    6  = 7.1225
    7  = 7.1194
    8  = 7.0972
    9  = 7.1132
    10 = 7.1169
    11 = 7.0858
    12 = 7.1279
    13 = 7.1324
}

foreach ($row in $caiLat.Keys) {
    $cell = $wsCAI.Cells.Item($row, 2)
    $cell.Value = $caiLat[$row]
    $cell.NumberFormat = "General"
}

$wsCAI.Cells.Item(4, 3).Value = -73.133

$wsCAI.Range("E7").Select()

$wsROBOS = $wb.Worksheets.Item("ROBOS")

$robosLat = @{
    2  = 7.15135
    3  = 7.12
    4  = 7.17055
    5  = 7.18668
    6  = 7.17804
    7  = 7.11633
    8  = 7.116
    9  = 7.091882
    10 = 7.117348
    11 = 7.118031
    12 = 7.125948
    13 = 7.128164
    14 = 7.14255
    15 = 7.149916
    16 = 7.133093
    17 = 7.15135
    18 = 7.12
    19 = 7.17055
    20 = 7.18668
    21 = 7.17804
    22 = 7.11633
    23 = 7.116
    24 = 7.091882
    25 = 7.117348
    26 = 7.118031
    27 = 7.125948
    28 = 7.128164
    29 = 7.14255
    30 = 7.149916
    31 = 7.133093
}

foreach ($row in $robosLat.Keys) {
    $cell = $wsROBOS.Cells.Item($row, 2)
    $cell.Value = $robosLat[$row]
    $cell.NumberFormat = "General"
}

$wsROBOS.Range("F23").Select()
